$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text for column G (row 2)
$ws.Range("G2").Value = "Other Peak(s) [Give a list of comma seperated list of rts and conc values, Example: ((2.11, 0.13), (3.25, 0.25))"

# G3 shares the old text of G4; update to new values per diff
$ws.Range("G3").Value = "(2.11, 0.13), (3.25, 0.25)"
$ws.Range("G4").Value = "(3.46, 2.99)"

# G5 was blank, now contains new text
$ws.Range("G5").Value = "string to fail here"

# G6, G7, G8 values cleared (were numeric, now blank)
$ws.Range("G6").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("G8").ClearContents()

# C7 and C8 were blank, now contain "ug/ml"
$ws.Range("C7").Value = "ug/ml"
$ws.Range("C8").Value = "ug/ml"

# Delete rows 13 through 23 (removes trailing empty styled rows)
$ws.Range("A13:G23").EntireRow.Delete()

# Widen column G (target stored OOXML width 106.98; this runtime stores
# width = ColumnWidth + 0.8333 snapped to 1/6 increments, so 106.15 is the
# input value that round-trips to the closest achievable stored width)
$ws.Columns.Item(7).ColumnWidth = 106.15

# Update selected cell
$ws.Range("C6").Select()

$wb.Save()
